# Apply data-cleaning fixes to the VERMONT_2015 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns to clean machine-readable names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Normalize "de" -> "De" in specific municipality / state names
$ws.Range("B4").Value = "Comitán De Domínguez"
$ws.Range("B8").Value = "Salto De Agua"
$ws.Range("A10").Value = "Ciudad De México"
$ws.Range("A15").Value = "Estado De México"
$ws.Range("B16").Value = "Ecatepec De Morelos"
$ws.Range("B23").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B28").Value = "Ixtlán De Juárez"
$ws.Range("B31").Value = "Tlalixtac De Cabrera"

# 3. Remove trailing metadata/footnote rows (67-71), keeping data through row 65
$ws.Range("A67:A71").EntireRow.Delete()
